$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking price values so Excel
# does not silently coerce them into floating point numbers.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.290.87"
$ws.Range("E2").Value = "  +0.20%  "

$ws.Range("D3").Value = "1.858.56"
$ws.Range("E3").Value = "  -0.06%  "

$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("E5").Value = "  +0.38%  "

$ws.Range("D6").Value = "238.27"
$ws.Range("E6").Value = "  +0.36%  "

$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "0.07886"
$ws.Range("E8").Value = "  +2.84%  "

$ws.Range("E9").Value = "  -0.45%  "

$ws.Range("D10").Value = "24.34"
$ws.Range("E10").Value = "  +4.36%  "

$ws.Range("D11").Value = "0.08179"
$ws.Range("E11").Value = "  +0.31%  "

$ws.Range("D12").Value = "1.864.15"
$ws.Range("E12").Value = "  -0.12%  "

$ws.Range("D13").Value = "0.7224"
$ws.Range("E13").Value = "  +0.62%  "

$ws.Range("D14").Value = "5.211"
$ws.Range("E14").Value = "  +1.10%  "

$ws.Range("D15").Value = "89.61"
$ws.Range("E15").Value = "  +0.05%  "

$ws.Range("D16").Value = "29.337.37"
$ws.Range("E16").Value = "  +0.29%  "

$ws.Range("D17").Value = "5.813"
$ws.Range("E17").Value = "  +0.98%  "

$ws.Range("D18").Value = "0.000007819"
$ws.Range("E18").Value = "  +1.51%  "

$ws.Range("D19").Value = "13.26"
$ws.Range("E19").Value = "  -0.16%  "

$ws.Range("D20").Value = "238.20"
$ws.Range("E20").Value = "  +0.17%  "

$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.121.70"
$ws.Range("E21").Value = "  +0.16%  "

$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.14%  "

$ws.Range("E23").Value = "  -0.12%  "

$ws.Range("D24").Value = "7.556"
$ws.Range("E24").Value = "  +1.12%  "

$ws.Range("D25").Value = "161.94"
$ws.Range("E25").Value = "  -0.54%  "

$ws.Range("D26").Value = "8.903"
$ws.Range("E26").Value = "  -1.25%  "

$ws.Range("D27").Value = "0.1423"
$ws.Range("E27").Value = "  -3.43%  "

$ws.Range("D28").Value = "18.10"
$ws.Range("E28").Value = "  +0.30%  "

$ws.Range("D29").Value = "1.915"
$ws.Range("E29").Value = "  -4.64%  "

$ws.Range("D30").Value = "1.386"
$ws.Range("E30").Value = "  -2.36%  "

$ws.Range("D31").Value = "1.475"
$ws.Range("E31").Value = "  -0.47%  "

$ws.Range("D32").Value = "4.318"
$ws.Range("E32").Value = "  -2.60%  "

$ws.Range("D33").Value = "4.055"
$ws.Range("E33").Value = "  +1.13%  "

$ws.Range("D34").Value = "0.05182"
$ws.Range("E34").Value = "  -0.41%  "

$ws.Range("D35").Value = "1.174"
$ws.Range("E35").Value = "  +0.65%  "

$ws.Range("D36").Value = "0.7137"
$ws.Range("E36").Value = "  +0.39%  "

$ws.Range("D37").Value = "0.9995"
$ws.Range("E37").Value = "  -0.16%  "

$ws.Range("D38").Value = "2.678"
$ws.Range("E38").Value = "  +0.64%  "

$ws.Range("E39").Value = "  -0.30%  "

$ws.Range("D40").Value = "2.689"
$ws.Range("E40").Value = "  -1.34%  "

$ws.Range("D41").Value = "1.151.25"
$ws.Range("E41").Value = "  +0.20%  "

$ws.Range("D42").Value = "0.9230"
$ws.Range("E42").Value = "  -1.31%  "

$ws.Range("D43").Value = "5.958"
$ws.Range("E43").Value = "  +1.48%  "

$ws.Range("E44").Value = "  -0.67%  "

$ws.Range("D45").Value = "70.70"
$ws.Range("E45").Value = "  -0.49%  "

$ws.Range("E46").Value = "  +0.00%  "

$ws.Range("D47").Value = "101.58"
$ws.Range("E47").Value = "  -1.84%  "

$ws.Range("E48").Value = "  -2.96%  "

$ws.Range("D50").Value = "9.192"
$ws.Range("E50").Value = "  +0.37%  "

$ws.Range("D51").Value = "6.999"
$ws.Range("E51").Value = "  +0.33%  "
